# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting newly generated data as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1419
$ws1.Range("F7").Value  = 11932
$ws1.Range("F13").Value = 2570
$ws1.Range("F15").Value = 169
$ws1.Range("F17").Value = 5167
$ws1.Range("F21").Value = 11392
$ws1.Range("F22").Value = 11389
$ws1.Range("F27").Value = 53

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1419
$ws4.Range("F7").Value  = 11932
$ws4.Range("F13").Value = 2570
$ws4.Range("F14").Value = 4
$ws4.Range("F16").Value = 169
$ws4.Range("F18").Value = 5167
$ws4.Range("F22").Value = 11392
$ws4.Range("F23").Value = 11389
$ws4.Range("F28").Value = 53
